# CORELIMS-98 - production addition of F3
# Renumber the DRW drawer barcodes in column C: each barcode value "DRW<n>"
# becomes "DRW<n+49>". Each barcode value spans a block of 4 consecutive
# rows (one row per FB1..FB4 sub-location), starting at row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 2
$rowsPerGroup = 4
$groupCount = 42
$offset = 49
$startNumber = 520

for ($i = 0; $i -lt $groupCount; $i++) {
    $row = $startRow + ($i * $rowsPerGroup)
    $newNumber = $startNumber + $i + $offset
    $newValue = "DRW" + $newNumber
    $endRow = $row + $rowsPerGroup - 1
    $rangeAddress = "C" + $row + ":C" + $endRow
    $ws.Range($rangeAddress).Value = $newValue
}
